$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DATE values (column A) -- entered as text via formula+paste-special so they
# stay as shared-string text instead of being auto-converted to date serials,
# matching the source data convention used throughout this column.
$dates = @("2021-02-25", "2021-02-26", "2021-02-27", "2021-02-28", "2021-03-01", "2021-03-02", "2021-03-03", "2021-03-04", "2021-03-05")

$startRow = 367
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = '="' + $dates[$i] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# Numeric data for columns B:O
$rowsData = @(
    @(2180, 151, 1327, 702, 2784, 79, 482, 2223, 6.93, 60.87, 32.2, 2.84, 17.31, 79.84999999999999),
    @(2182, 138, 1293, 751, 2784, 73, 484, 2227, 6.32, 59.26, 34.42, 2.62, 17.39, 79.98999999999999),
    @(2182, 116, 1202, 864, 2785, 67, 505, 2213, 5.32, 55.09, 39.6, 2.41, 18.13, 79.45999999999999),
    @(2175, 124, 1175, 876, 2786, 66, 490, 2230, 5.7, 54.02, 40.28, 2.37, 17.59, 80.04000000000001),
    @(2182, 123, 1176, 883, 2787, 64, 501, 2222, 5.64, 53.9, 40.47, 2.3, 17.98, 79.73),
    @(2186, 131, 1311, 744, 2786, 65, 500, 2221, 5.99, 59.97, 34.03, 2.33, 17.95, 79.72),
    @(2183, 131, 1343, 709, 2785, 63, 512, 2210, 6, 61.52, 32.48, 2.26, 18.38, 79.34999999999999),
    @(2183, 123, 1314, 746, 2778, 64, 484, 2230, 5.63, 60.19, 34.17, 2.3, 17.42, 80.27),
    @(2180, 118, 1340, 722, 2783, 60, 487, 2236, 5.41, 61.47, 33.12, 2.16, 17.5, 80.34)
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $row = $startRow + $i
    $vals = $rowsData[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
